$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = "Iedereen heeft goed gewerkt"

$ws.Range("A7").Select()
